$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.56"
$ws.Range("E2").Value = "'-0.73%"

$ws.Range("D3").Value = "'31.45"
$ws.Range("E3").Value = "'-0.07%"

$ws.Range("D4").Value = "'5.094"
$ws.Range("E4").Value = "'-1.01%"

$ws.Range("D5").Value = "'0.08005"
$ws.Range("E5").Value = "'9.37%"

$ws.Range("D6").Value = "'2.451"
$ws.Range("E6").Value = "'34.39%"

$ws.Range("D7").Value = "'7.789"
$ws.Range("E7").Value = "'0.29%"

$ws.Range("D8").Value = "'0.9227"
$ws.Range("E8").Value = "'-0.25%"

$ws.Range("D9").Value = "'0.1756"
$ws.Range("E9").Value = "'4.80%"

$ws.Range("D10").Value = "'0.07341"
$ws.Range("E10").Value = "'2.91%"

$ws.Range("D11").Value = "'0.08812"
$ws.Range("E11").Value = "'9.03%"

$ws.Range("E12").Value = "'1.06%"

$ws.Range("D13").Value = "'0.09996"
$ws.Range("E13").Value = "'0.88%"

$ws.Range("D14").Value = "'0.001496"
$ws.Range("E14").Value = "'0.65%"

$ws.Range("D15").Value = "'0.006008"
$ws.Range("E15").Value = "'-2.63%"

$ws.Range("D16").Value = "'3.505"
$ws.Range("E16").Value = "'1.43%"

$ws.Range("D17").Value = "'3.802"
$ws.Range("E17").Value = "'1.67%"

$ws.Range("D18").Value = "'2.248"
$ws.Range("E18").Value = "'1.17%"

$ws.Range("E19").Value = "'1.86%"

$ws.Range("E20").Value = "'1.48%"

$ws.Range("D21").Value = "'4.291"
$ws.Range("E21").Value = "'-5.83%"

$ws.Range("D22").Value = "'0.1618"
$ws.Range("E22").Value = "'2.29%"

$ws.Range("D23").Value = "'0.04597"
$ws.Range("E23").Value = "'-1.02%"

$ws.Range("D24").Value = "'0.001243"
$ws.Range("E24").Value = "'2.41%"

$ws.Range("D25").Value = "'0.004428"
$ws.Range("E25").Value = "'-6.52%"

$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'-7.50%"

$ws.Range("D27").Value = "'0.0003429"
$ws.Range("E27").Value = "'83.09%"

$ws.Range("D39").Value = "'0.01774"
$ws.Range("E39").Value = "'3.23%"

$ws.Range("D40").Value = "'0.04461"
$ws.Range("E40").Value = "'-0.29%"

$ws.Range("D41").Value = "'0.006966"
$ws.Range("E41").Value = "'-1.52%"

$ws.Range("D42").Value = "'0.1344"
$ws.Range("E42").Value = "'0.80%"

$ws.Range("D43").Value = "'0.002211"
$ws.Range("E43").Value = "'3.97%"

$ws.Range("E44").Value = "'-6.17%"

$ws.Range("D45").Value = "'0.00006562"
$ws.Range("E45").Value = "'5.47%"

$ws.Range("E46").Value = "'0.05%"

$ws.Range("E48").Value = "'-55.55%"

$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.05%"

$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.12%"
